$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (row 5) - the new TPM data only has 3 data rows.
$ws.Rows(5).Delete()

# Row 2: FAPs -> Rbp4 -> Stra6 -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Rbp4"
$ws.Range("C2").Value = "Stra6"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.656305333333333
$ws.Range("H2").Value = 4.968916
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.033075
$ws.Range("N2").Value = 0.09922499999999999
$ws.Range("O2").Value = 0.01056319585618255
$ws.Range("P2").Value = 0.01056319585618256
$ws.Range("Q2").Value = 0.0547822989
$ws.Range("R2").Value = 0.4930406901
$ws.Range("S2").Value = 0.01056319585618255
$ws.Range("T2").Value = 0.01056319585618256

# Row 3: FAPs -> Rbp4 -> Stra6 -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Rbp4"
$ws.Range("C3").Value = "Stra6"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.656305333333333
$ws.Range("H3").Value = 4.968916
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 1.307784333333333
$ws.Range("N3").Value = 3.923353
$ws.Range("O3").Value = 0.4176683915539571
$ws.Range("P3").Value = 0.4176683915539571
$ws.Range("Q3").Value = 2.166090166149778
$ws.Range("R3").Value = 19.494811495348
$ws.Range("S3").Value = 0.4176683915539571
$ws.Range("T3").Value = 0.4176683915539571

# Row 4: FAPs -> Rbp4 -> Stra6 -> MuSCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Rbp4"
$ws.Range("C4").Value = "Stra6"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.656305333333333
$ws.Range("H4").Value = 4.968916
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.790295333333333
$ws.Range("N4").Value = 5.370886
$ws.Range("O4").Value = 0.5717684125898603
$ws.Range("P4").Value = 0.5717684125898604
$ws.Range("Q4").Value = 2.965275708841777
$ws.Range("R4").Value = 26.687481379576
$ws.Range("S4").Value = 0.5717684125898603
$ws.Range("T4").Value = 0.5717684125898604

Write-Output "done"
